$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 4 corresponds to 6525353b-4c16-4cfb-be57-d26d14a0f9a4...
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-29 03:50:52"
$wsZhCn.Range("G4").Value = "2016-02-29 03:51:40"

# de-de sheet: row 4 corresponds to 6525353b-4c16-4cfb-be57-d26d14a0f9a4...
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-29 03:51:03"
$wsDeDe.Range("G4").Value = "2016-02-29 03:52:00"
